$wb = $excel.ActiveWorkbook

$wsARM = $wb.Worksheets.Item("ARM")
$wsCUL = $wb.Worksheets.Item("CUL")

# Row 13
$wsCUL.Range("H13").Value = 1911.1111
$wsCUL.Range("I13").Value = 300
$wsCUL.Range("J13").Value = 2112.5
$wsCUL.Range("K13").Value = 900
$wsCUL.Range("L13").Value = 6337.5
$wsCUL.Range("M13").Value = -732
$wsCUL.Range("N13").Value = -6673.5

# Row 62
$wsCUL.Range("H62").Value = 1683.1578
$wsCUL.Range("I62").Value = 1763.5294
$wsCUL.Range("J62").Value = 1000
$wsCUL.Range("K62").Value = 5290.5882
$wsCUL.Range("L62").Value = 3000
$wsCUL.Range("M62").Value = -4604.5882
$wsCUL.Range("N62").Value = -4372

# Row 63
$wsCUL.Range("H63").Value = 14755.556
$wsCUL.Range("I63").Value = 10000
$wsCUL.Range("J63").Value = 15350
$wsCUL.Range("K63").Value = 30000
$wsCUL.Range("L63").Value = 46050
$wsCUL.Range("M63").Value = -29251
$wsCUL.Range("N63").Value = -47548

# Row 64
$wsCUL.Range("H64").Value = 3200
$wsCUL.Range("I64").Value = 2666.6667
$wsCUL.Range("K64").Value = 8000.000100000001
$wsCUL.Range("M64").Value = -7730.000100000001

# Row 65
$wsCUL.Range("H65").Value = 1683.1578
$wsCUL.Range("I65").Value = 1763.5294
$wsCUL.Range("J65").Value = 1000
$wsCUL.Range("K65").Value = 15871.7646
$wsCUL.Range("L65").Value = 9000
$wsCUL.Range("M65").Value = -12439.7646
$wsCUL.Range("N65").Value = -15864

# Row 66
$wsCUL.Range("H66").Value = 14755.556
$wsCUL.Range("I66").Value = 10000
$wsCUL.Range("J66").Value = 15350
$wsCUL.Range("K66").Value = 90000
$wsCUL.Range("L66").Value = 138150
$wsCUL.Range("M66").Value = -86256
$wsCUL.Range("N66").Value = -145638

# Row 67
$wsCUL.Range("H67").Value = 3200
$wsCUL.Range("I67").Value = 2666.6667
$wsCUL.Range("K67").Value = 8000.000100000001
$wsCUL.Range("M67").Value = -7064.000100000001

# Row 68
$wsCUL.Range("H68").Value = 1211.5714
$wsCUL.Range("I68").Value = 427.33334
$wsCUL.Range("J68").Value = 1799.75
$wsCUL.Range("K68").Value = 1282.00002
$wsCUL.Range("L68").Value = 5399.25
$wsCUL.Range("M68").Value = -471.0000199999999
$wsCUL.Range("N68").Value = -7021.25

# Row 69
$wsCUL.Range("H69").Value = 1028
$wsCUL.Range("I69").Value = 808
$wsCUL.Range("J69").Value = 1248
$wsCUL.Range("K69").Value = 2424
$wsCUL.Range("L69").Value = 3744
$wsCUL.Range("M69").Value = -1613
$wsCUL.Range("N69").Value = -5366

# Row 70
$wsCUL.Range("H70").Value = 2666.6667

# Row 71
$wsCUL.Range("H71").Value = 1211.5714
$wsCUL.Range("I71").Value = 427.33334
$wsCUL.Range("J71").Value = 1799.75
$wsCUL.Range("K71").Value = 3846.00006
$wsCUL.Range("L71").Value = 16197.75
$wsCUL.Range("M71").Value = 209.9999399999997
$wsCUL.Range("N71").Value = -24309.75

# Row 72
$wsCUL.Range("H72").Value = 1028
$wsCUL.Range("I72").Value = 808
$wsCUL.Range("J72").Value = 1248
$wsCUL.Range("K72").Value = 7272
$wsCUL.Range("L72").Value = 11232
$wsCUL.Range("M72").Value = -3216
$wsCUL.Range("N72").Value = -19344

# Row 73
$wsCUL.Range("H73").Value = 2666.6667

# Row 74
$wsCUL.Range("H74").Value = 8752.166999999999
$wsCUL.Range("I74").Value = 5256.5
$wsCUL.Range("J74").Value = 10500
$wsCUL.Range("K74").Value = 15769.5
$wsCUL.Range("L74").Value = 31500
$wsCUL.Range("M74").Value = -14708.5
$wsCUL.Range("N74").Value = -33622

# Row 75
$wsCUL.Range("H75").Value = 2922.4
$wsCUL.Range("I75").Value = 1400
$wsCUL.Range("J75").Value = 3091.5557
$wsCUL.Range("K75").Value = 4200
$wsCUL.Range("L75").Value = 9274.667099999999
$wsCUL.Range("M75").Value = -3202
$wsCUL.Range("N75").Value = -11270.6671

# Row 76
$wsCUL.Range("H76").Value = 4000
$wsCUL.Range("I76").Value = 0
$wsCUL.Range("J76").Value = 4000
$wsCUL.Range("K76").Value = 0
$wsCUL.Range("L76").Value = 12000
$wsCUL.Range("M76").ClearContents()
$wsCUL.Range("N76").Value = -12766

# Row 77
$wsCUL.Range("H77").Value = 8752.166999999999
$wsCUL.Range("I77").Value = 5256.5
$wsCUL.Range("J77").Value = 10500
$wsCUL.Range("K77").Value = 47308.5
$wsCUL.Range("L77").Value = 94500
$wsCUL.Range("M77").Value = -42004.5
$wsCUL.Range("N77").Value = -105108

# Row 78
$wsCUL.Range("H78").Value = 2922.4
$wsCUL.Range("I78").Value = 1400
$wsCUL.Range("J78").Value = 3091.5557
$wsCUL.Range("K78").Value = 12600
$wsCUL.Range("L78").Value = 27824.0013
$wsCUL.Range("M78").Value = -7608
$wsCUL.Range("N78").Value = -37808.0013

# Row 79
$wsCUL.Range("H79").Value = 4000
$wsCUL.Range("I79").Value = 0
$wsCUL.Range("J79").Value = 4000
$wsCUL.Range("K79").Value = 0
$wsCUL.Range("L79").Value = 12000
$wsCUL.Range("M79").ClearContents()
$wsCUL.Range("N79").Value = -14652

# Row 80
$wsCUL.Range("H80").Value = 5750.0625
$wsCUL.Range("I80").Value = 2666.6667
$wsCUL.Range("J80").Value = 6461.615
$wsCUL.Range("K80").Value = 8000.000100000001
$wsCUL.Range("L80").Value = 19384.845
$wsCUL.Range("M80").Value = -7064.000100000001
$wsCUL.Range("N80").Value = -21256.845

# Row 81
$wsCUL.Range("H81").Value = 2750
$wsCUL.Range("I81").Value = 1500
$wsCUL.Range("K81").Value = 4500
$wsCUL.Range("M81").Value = -3377

# Row 82
$wsCUL.Range("H82").Value = 7201.625
$wsCUL.Range("I82").Value = 7156.5
$wsCUL.Range("J82").Value = 7216.6665
$wsCUL.Range("K82").Value = 21469.5
$wsCUL.Range("L82").Value = 21649.9995
$wsCUL.Range("M82").Value = -21063.5
$wsCUL.Range("N82").Value = -22461.9995

# Row 83
$wsCUL.Range("H83").Value = 5750.0625
$wsCUL.Range("I83").Value = 2666.6667
$wsCUL.Range("J83").Value = 6461.615
$wsCUL.Range("K83").Value = 24000.0003
$wsCUL.Range("L83").Value = 58154.535
$wsCUL.Range("M83").Value = -19320.0003
$wsCUL.Range("N83").Value = -67514.535

# Row 84
$wsCUL.Range("H84").Value = 2750
$wsCUL.Range("I84").Value = 1500
$wsCUL.Range("K84").Value = 13500
$wsCUL.Range("M84").Value = -7884

# Row 85
$wsCUL.Range("H85").Value = 7201.625
$wsCUL.Range("I85").Value = 7156.5
$wsCUL.Range("J85").Value = 7216.6665
$wsCUL.Range("K85").Value = 21469.5
$wsCUL.Range("L85").Value = 21649.9995
$wsCUL.Range("M85").Value = -20065.5
$wsCUL.Range("N85").Value = -24457.9995

# Row 86
$wsCUL.Range("H86").Value = 347.33334
$wsCUL.Range("I86").Value = 226
$wsCUL.Range("J86").Value = 387.77777
$wsCUL.Range("K86").Value = 678
$wsCUL.Range("L86").Value = 1163.33331
$wsCUL.Range("M86").Value = 508
$wsCUL.Range("N86").Value = -3535.33331

# Row 87
$wsCUL.Range("H87").Value = 1014
$wsCUL.Range("I87").Value = 1014
$wsCUL.Range("K87").Value = 3042
$wsCUL.Range("M87").Value = -1794

# Row 88
$wsCUL.Range("H88").Value = 3333.3333
$wsCUL.Range("J88").Value = 3333.3333
$wsCUL.Range("L88").Value = 9999.999899999999
$wsCUL.Range("N88").Value = -10855.9999

# Row 89
$wsCUL.Range("H89").Value = 347.33334
$wsCUL.Range("I89").Value = 226
$wsCUL.Range("J89").Value = 387.77777
$wsCUL.Range("K89").Value = 2034
$wsCUL.Range("L89").Value = 3489.99993
$wsCUL.Range("M89").Value = 3894
$wsCUL.Range("N89").Value = -15345.99993

# Row 90
$wsCUL.Range("H90").Value = 1014
$wsCUL.Range("I90").Value = 1014
$wsCUL.Range("K90").Value = 9126
$wsCUL.Range("M90").Value = -2886

# Row 91
$wsCUL.Range("H91").Value = 3333.3333
$wsCUL.Range("J91").Value = 3333.3333
$wsCUL.Range("L91").Value = 9999.999899999999
$wsCUL.Range("N91").Value = -12963.9999

# Row 102
$wsARM.Range("H102").Value = 76925260
$wsARM.Range("I102").Value = 100001940
$wsARM.Range("J102").Value = 3007.3333
$wsARM.Range("K102").Value = 100001940
$wsARM.Range("L102").Value = 3007.3333
$wsARM.Range("M102").Value = -100000318
$wsARM.Range("N102").Value = -6251.3333

# Row 131
$wsCUL.Range("H131").Value = 899.35
$wsCUL.Range("I131").Value = 0
$wsCUL.Range("J131").Value = 899.35
$wsCUL.Range("K131").Value = 0
$wsCUL.Range("L131").Value = 2698.05
$wsCUL.Range("M131").ClearContents()
$wsCUL.Range("N131").Value = -12778.05

# Row 134
$wsCUL.Range("H134").Value = 3307.9614
$wsCUL.Range("I134").Value = 2837.2104
$wsCUL.Range("J134").Value = 4585.7144
$wsCUL.Range("K134").Value = 8511.6312
$wsCUL.Range("L134").Value = 13757.1432
$wsCUL.Range("M134").Value = -3441.6312
$wsCUL.Range("N134").Value = -23897.1432
